# Node-RED template update: "update nodered node to support connect with uns"
#
# The header row (row 1) drops the "FileName(Required)" column (old column B),
# renames "FolderPath(" -> "FilePath(" in column A, and tweaks the type-list
# comment from "int" to "integer". Deleting the column shifts the remaining
# columns (FileAlias, AttributeName, AttributeType, TagConfiguration) one slot
# to the left, so their cell comments have to be re-anchored to the new
# addresses (comment anchors do not automatically follow a column delete).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Capture the comment text of the columns that will shift left (D, E, F -> C, D, E)
# before touching the sheet, since comment anchors stay put on column delete.
$textAttrName  = $ws.Range("D1").Comment.Text()
$textAttrType  = $ws.Range("E1").Comment.Text()
$textTagConfig = $ws.Range("F1").Comment.Text()

# A1: "FolderPath(Required)" -> "FilePath(Required)" (keep the bold "Required)" run).
$ws.Range("A1").Characters(1, 11).Text = "FilePath("

# Drop every existing comment on row 1; the ones that survive the column
# delete get re-created below at their correct, post-delete addresses.
$ws.Range("B1").Comment.Delete()
$ws.Range("D1").Comment.Delete()
$ws.Range("E1").Comment.Delete()
$ws.Range("F1").Comment.Delete()

# Delete column B entirely (the "FileName(Required)" column).
# C,D,E,F ("FileAlias", "AttributeName(Required)", "AttributeType(Required)",
# "TagConfiguration(Required)") shift left to B,C,D,E.
$ws.Columns("B").Delete()

# Re-create the surviving comments on their new cell addresses. The
# AttributeType comment's enumerated list changes "int" to "integer".
$ws.Range("C1").AddComment($textAttrName)
$ws.Range("D1").AddComment(($textAttrType -replace '\bint\b', 'integer'))
$ws.Range("E1").AddComment($textTagConfig)

# Match the target view state (active cell/selection).
$ws.Range("E7").Select()
